# Updated symbol list on Mon Dec 12 20:05:57 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) updates for specific rows
$ws.Range("D2").Value  = "274.84"
$ws.Range("D3").Value  = "21.13"
$ws.Range("D4").Value  = "6.215"
$ws.Range("D5").Value  = "0.06179"
$ws.Range("D6").Value  = "3.576"
$ws.Range("D7").Value  = "1.529"
$ws.Range("D8").Value  = "6.527"
$ws.Range("D9").Value  = "0.8233"
$ws.Range("D10").Value = "0.1648"
$ws.Range("D11").Value = "0.08211"
$ws.Range("D12").Value = "0.03431"
$ws.Range("D13").Value = "0.03125"
$ws.Range("D14").Value = "0.09132"
$ws.Range("D15").Value = "3.769"
$ws.Range("D16").Value = "0.001618"
$ws.Range("D17").Value = "0.04692"
$ws.Range("D18").Value = "0.006430"
$ws.Range("D19").Value = "0.006135"
$ws.Range("D21").Value = "0.0001500"
$ws.Range("D22").Value = "3.722"
$ws.Range("D25").Value = "0.3324"
$ws.Range("D26").Value = "0.1232"
$ws.Range("D40").Value = "0.04743"
$ws.Range("D41").Value = "0.005500"
$ws.Range("D42").Value = "0.007028"
$ws.Range("D43").Value = "0.1107"
$ws.Range("D44").Value = "0.01035"
$ws.Range("D45").Value = "0.00006275"
$ws.Range("D47").Value = "0.7230"
$ws.Range("D49").Value = "0.00001900"

# Column E (Volume(1h)) label updates
$ws.Range("E41").Value = "40CEJICEJI"
$ws.Range("E42").Value = "41KickTokenKICKBestin24h"

# Column G (Hora) updates for all data rows 2-51
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 7).Value = "20"
}
